# Create Session for client and save in database
# -> add an "english" helper column and a batch of new word-pairs to the
#    word-guessing list on Sheet1, then re-apply the existing sort/CF setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D: tag the first two existing rows as "english" --------
$ws.Cells.Item(1,4).Value = "english"
$ws.Cells.Item(2,4).Value = "english"

# --- Row numbers (column A) for the new word pairs ----------------------
$ws.Cells.Item(99,1).Value  = 99
$ws.Cells.Item(100,1).Value = 100
$ws.Cells.Item(101,1).Value = 101
$ws.Cells.Item(102,1).Value = 102
$ws.Cells.Item(103,1).Value = 103
$ws.Cells.Item(104,1).Value = 104
$ws.Cells.Item(105,1).Value = 105
$ws.Cells.Item(106,1).Value = 106
$ws.Cells.Item(107,1).Value = 107
$ws.Cells.Item(108,1).Value = 108
$ws.Cells.Item(109,1).Value = 109
$ws.Cells.Item(110,1).Value = 110

# --- New word pairs (columns B/C), entered in the same order the author
#     originally typed them so shared-string ids line up -----------------
$ws.Cells.Item(99,2).Value  = "甄子丹"
$ws.Cells.Item(99,3).Value  = "李连杰"

$ws.Cells.Item(100,2).Value = "狼"
$ws.Cells.Item(100,3).Value = "狐狸"

$ws.Cells.Item(101,2).Value = "成龙"
$ws.Cells.Item(101,3).Value = "李连杰"

$ws.Cells.Item(102,2).Value = "郭靖"
$ws.Cells.Item(102,3).Value = "杨过"

$ws.Cells.Item(103,2).Value = "红孩儿"
$ws.Cells.Item(103,3).Value = "哪吒"

$ws.Cells.Item(105,2).Value = "Punggol"
$ws.Cells.Item(105,3).Value = "Seng Kang"

$ws.Cells.Item(106,3).Value = "NTU"
$ws.Cells.Item(106,2).Value = "NUS"

$ws.Cells.Item(104,3).Value = "Buger King"

$ws.Cells.Item(107,2).Value = "Badminton"
$ws.Cells.Item(107,3).Value = "Tennis"

$ws.Cells.Item(104,2).Value = "McDonalds"

$ws.Cells.Item(108,2).Value = "白骨精"
$ws.Cells.Item(108,3).Value = "蜘蛛精"

$ws.Cells.Item(109,2).Value = "Winnie the Pooh"
$ws.Cells.Item(109,3).Value = "Piglet"

$ws.Cells.Item(110,2).Value = "樱木花道"
$ws.Cells.Item(110,3).Value = "流川枫"

# --- Re-apply the sheet's sort (keeps column A ascending, but refreshes
#     the stored sortState reference to cover the newly added rows) ------
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A1:A125"))
$sortObj.SetRange($ws.Range("A1:C125"))
$sortObj.Header = 0
$sortObj.Apply()

# --- Restore first/last priority ordering on the duplicate-value rules --
$cfRange = $ws.Range("B1:C1048576")
$fcs = $cfRange.FormatConditions
$fcs.Item(1).SetFirstPriority()

# --- Update the view: scroll down and select the last cell touched ------
$ws.Range("A95").Select()
$ws.Range("J109").Select()
